# Weekly update: insert a new price record for "Poroto granado" at row 20,
# pushing the existing rows (old row 20 .. old row 54) down by one row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 20 - shifts rows 20..54 down to 21..55
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the new record's data
$ws.Cells.Item(20, 1).Value = 11
$ws.Cells.Item(20, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(20, 3).Value = "Bíobío"
$ws.Cells.Item(20, 4).Value = 44952
$ws.Cells.Item(20, 5).Value = 8
$ws.Cells.Item(20, 6).Value = 100112030
$ws.Cells.Item(20, 7).Value = "Poroto granado"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 38000
$ws.Cells.Item(20, 12).Value = 40000
$ws.Cells.Item(20, 13).Value = 39000
$ws.Cells.Item(20, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(20, 15).Value = "Región Metropolitana"
$ws.Cells.Item(20, 16).Value = 1560
$ws.Cells.Item(20, 17).Value = 25
$ws.Cells.Item(20, 18).Value = "Hortaliza"
